$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.6753301551942219
$ws.Range("C2").Value = 0.3127903958511391
$ws.Range("D2").Value = 0.1575252929769615
$ws.Range("E2").Value = 0.496779210170732
$ws.Range("G2").Value = 1.642425054193055

# Row 3
$ws.Range("B3").Value = 0.6753301551942219
$ws.Range("C3").Value = 0.3127903958511391
$ws.Range("D3").Value = 3.900430680208489
$ws.Range("E3").Value = 8.660232485948974
$ws.Range("G3").Value = 13.54878371720282
